$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 62, shifting existing rows 62..164 down to 63..165
$ws.Rows.Item(62).Insert()

# Populate the newly inserted row 62 with the new data record
$ws.Cells.Item(62, 1).Value = 4
$ws.Cells.Item(62, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(62, 3).Value = "Los Lagos"
$ws.Cells.Item(62, 4).Value = 44581
$ws.Cells.Item(62, 5).Value = 10
$ws.Cells.Item(62, 6).Value = 100112039
$ws.Cells.Item(62, 7).Value = "Ciboulette"
$ws.Cells.Item(62, 8).Value = "Sin especificar"
$ws.Cells.Item(62, 9).Value = "Primera"
$ws.Cells.Item(62, 10).Value = 80
$ws.Cells.Item(62, 11).Value = 2500
$ws.Cells.Item(62, 12).Value = 2500
$ws.Cells.Item(62, 13).Value = 2500
$ws.Cells.Item(62, 14).Value = "`$/docena de atados"
$ws.Cells.Item(62, 15).Value = "Región Metropolitana"
$ws.Cells.Item(62, 16).Value = 833
$ws.Cells.Item(62, 17).Value = 3
$ws.Cells.Item(62, 18).Value = "Hortaliza"

# Match the date-cell style used by the rest of column D
$ws.Cells.Item(62, 4).NumberFormat = $ws.Cells.Item(63, 4).NumberFormat
